$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (RMSE), C (DM stat), D (Encompassing t-stat β)
# Rows 2-12 correspond to Horizon 1-11
$data = @(
    @(89.933845520019531, -0.82123315334320068, 0.24532508850097656),
    @(87.274131774902344, -1.0597113370895386, 0.91008186340332031),
    @(91.155021667480469, -0.89096373319625854, 2.4081621170043945),
    @(92.403213500976562, -1.1595971584320068, 1.2837705612182617),
    @(100.34432220458984, 0.10446625947952271, -0.71237343549728394),
    @(100.99150085449219, 0.17972996830940247, 1.9403207302093506),
    @(95.411331176757812, -0.65167880058288574, 1.968997597694397),
    @(99.457733154296875, -0.25193247199058533, 0.60457837581634521),
    @(115.80441284179688, 4.3055582046508789, 0.5911906361579895),
    @(122.04225921630859, 3.9517755508422852, 0.6602790355682373),
    @(106.31224822998047, 2.2163121700286865, -0.094798311591148376)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
}
